{"js": "// Replace the two-digit multiplication problems in the document's table\n// with their updated values, per the commit's diff. Each \"before\" string\n// is unique in the document, so a case-sensitive search for the exact\n// literal text (including the trailing \"=\") is safe and unambiguous.\nconst replacements = [\n  [\"52\u00d793=\", \"97\u00d773=\"],\n  [\"78\u00d711=\", \"72\u00d713=\"],\n  [\"52\u00d731=\", \"42\u00d772=\"],\n  [\"42\u00d782=\", \"88\u00d730=\"],\n  [\"46\u00d749=\", \"15\u00d789=\"],\n  [\"86\u00d763=\", \"76\u00d715=\"],\n  [\"93\u00d749=\", \"22\u00d771=\"],\n  [\"39\u00d762=\", \"11\u00d718=\"],\n  [\"52\u00d772=\", \"83\u00d716=\"],\n  [\"57\u00d798=\", \"93\u00d725=\"],\n  [\"19\u00d786=\", \"81\u00d760=\"],\n  [\"55\u00d784=\", \"66\u00d761=\"],\n  [\"66\u00d738=\", \"32\u00d756=\"],\n  [\"95\u00d781=\", \"66\u00d788=\"],\n  [\"48\u00d789=\", \"60\u00d791=\"],\n  [\"28\u00d715=\", \"74\u00d764=\"],\n  [\"46\u00d726=\", \"66\u00d771=\"],\n  [\"74\u00d712=\", \"77\u00d726=\"],\n  [\"65\u00d712=\", \"49\u00d752=\"],\n  [\"31\u00d722=\", \"60\u00d745=\"],\n  [\"52\u00d718=\", \"46\u00d746=\"],\n  [\"73\u00d754=\", \"58\u00d764=\"],\n  [\"63\u00d798=\", \"13\u00d771=\"],\n  [\"87\u00d754=\", \"93\u00d788=\"],\n  [\"71\u00d779=\", \"96\u00d724=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the document's table\n# with their updated values, per the commit's diff. Each \"before\" string is\n# unique in the document, so Find/Replace-All on the exact literal text is\n# safe and unambiguous. (Note: we read each pair's strings straight out of\n# $pairs/$p rather than reading them back off the Find/Replacement COM\n# properties, since re-reading those properties can mangle non-ASCII\n# characters such as \"\u00d7\".)\n$doc = $word.ActiveDocument\n\n$pairs = @(\n  @(\"52\u00d793=\", \"97\u00d773=\"),\n  @(\"78\u00d711=\", \"72\u00d713=\"),\n  @(\"52\u00d731=\", \"42\u00d772=\"),\n  @(\"42\u00d782=\", \"88\u00d730=\"),\n  @(\"46\u00d749=\", \"15\u00d789=\"),\n  @(\"86\u00d763=\", \"76\u00d715=\"),\n  @(\"93\u00d749=\", \"22\u00d771=\"),\n  @(\"39\u00d762=\", \"11\u00d718=\"),\n  @(\"52\u00d772=\", \"83\u00d716=\"),\n  @(\"57\u00d798=\", \"93\u00d725=\"),\n  @(\"19\u00d786=\", \"81\u00d760=\"),\n  @(\"55\u00d784=\", \"66\u00d761=\"),\n  @(\"66\u00d738=\", \"32\u00d756=\"),\n  @(\"95\u00d781=\", \"66\u00d788=\"),\n  @(\"48\u00d789=\", \"60\u00d791=\"),\n  @(\"28\u00d715=\", \"74\u00d764=\"),\n  @(\"46\u00d726=\", \"66\u00d771=\"),\n  @(\"74\u00d712=\", \"77\u00d726=\"),\n  @(\"65\u00d712=\", \"49\u00d752=\"),\n  @(\"31\u00d722=\", \"60\u00d745=\"),\n  @(\"52\u00d718=\", \"46\u00d746=\"),\n  @(\"73\u00d754=\", \"58\u00d764=\"),\n  @(\"63\u00d798=\", \"13\u00d771=\"),\n  @(\"87\u00d754=\", \"93\u00d788=\"),\n  @(\"71\u00d779=\", \"96\u00d724=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($p in $pairs) {\n  $before = $p[0]\n  $after = $p[1]\n\n  $find = $doc.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $before\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $after\n  $find.Execute($before, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $after, $wdReplaceAll)\n}\n"}
